$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Write_Review_Link_Xpath" value (B11)
#    "Write a review" -> "Write Your Review"
$ws.Range("B11").Value = "Write Your Review"

# B11's style changes from s="8" to s="9" (Courier New, size 10, blue FF2A00FF,
# same look as the existing xpath-value cells such as B29/B30). Copy the
# formatting from one of those cells instead of setting font properties one
# at a time so the workbook's style table stays clean.
$ws.Range("B29").Copy()
$ws.Range("B11").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# 2. Update the "Review_Date_Xpath" value (B32)
$ws.Range("B32").Value = "((//li[@class='rvw_title block clear']/div[2]))"

# 3. Append a new row 33 with a new locator entry
$ws.Range("A33").Value = "No_of_Customer_Reviews_DetailsPage_Xpath"
$ws.Range("B33").Value = "//p[@class='rating-links customreview']/a[1]"

# 4. Update the selected cell shown in the sheet view
$ws.Range("B11").Select()
